# POSTAL FINAL ERROR FIX
#
# Fills in the "REVISION" column (J) with the corrected postal code for every
# row that was previously left blank and could be resolved, and flags the
# handful of rows that could not be resolved (no corrected postal code
# available) with a yellow highlight across the row and a "?" placeholder in
# the REVISION column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Corrected postal codes (REVISION column, col J) for rows that had no
#        revision value yet. -------------------------------------------------
$revisions = @{
    152 = 61252
    153 = 41283
    159 = 17530
    161 = 17143
    163 = 17350
    164 = 16810
    167 = 61463
    168 = 40227
    169 = 40227
    171 = 24251
    172 = 27684
    173 = 28784
    174 = 28554
    176 = 41384
    177 = 22733
    178 = 17350
    180 = 80361
    181 = 28784
    182 = 20372
    183 = 28784
    184 = 20372
    185 = 28784
    186 = 20372
    187 = 20372
    188 = 20372
    189 = 28784
    190 = 20372
    191 = 20372
    192 = 20372
    194 = 20372
    195 = 28784
    196 = 29663
    198 = 21167
    200 = 26152
    202 = 20998
    203 = 22272
    204 = 28784
    205 = 29352
    206 = 59462
    208 = 57253
    210 = 41384
    211 = 28784
    212 = 28784
    213 = 20372
    214 = 28784
    217 = 20372
    218 = 28784
    219 = 20372
    220 = 28784
    223 = 28784
    224 = 80361
    226 = 17530
    230 = 17350
    232 = 20372
    233 = 57253
    236 = 28293
    238 = 28784
    239 = 20372
    241 = 30961
    243 = 20986
    245 = 61151
    246 = 16610
    247 = 15125
    251 = 17540
    252 = 17540
    256 = 17540
    257 = 21274
    260 = 28293
    261 = 28293
    267 = 21162
    270 = 42161
    274 = 22733
    275 = 34684
    276 = 20773
    278 = 28784
    280 = 21273
    286 = 37311
    293 = 35353
}

foreach ($row in $revisions.Keys) {
    $ws.Cells.Item($row, 10).Value = $revisions[$row]
}

# --- 2. Rows that could not be resolved: mark with "?" and highlight the
#        whole row in yellow so they stand out as still needing attention. --
$unresolvedRows = @(18, 23, 150, 193, 199)

foreach ($row in $unresolvedRows) {
    $ws.Cells.Item($row, 10).Value = "?"
    $rowRange = $ws.Range("A" + $row + ":J" + $row)
    $rowRange.Interior.Color = 65535
}

# --- 3. Leave the view roughly where the last edit happened (cosmetic). ----
$null = $ws.Range("E294").Select()
